$wb = $excel.ActiveWorkbook

$wsValues = $wb.Worksheets.Item("values")
$wsValues.Range("D2").Value = 0.03205053372796332
$wsValues.Range("E2").Value = -342529.9652282526

$wsCashflows = $wb.Worksheets.Item("Cashflows")
$wsCashflows.Range("D2").Value = 0.02508731863097187
$wsCashflows.Range("E2").Value = 128224.0730027451
$wsCashflows.Range("D3").Value = 0.03570358392086448
$wsCashflows.Range("E3").Value = 180501.4520443704
$wsCashflows.Range("D4").Value = 0.0329919142530939
$wsCashflows.Range("E4").Value = 168625.3395158133
$wsCashflows.Range("D5").Value = 0.03161977715018041
$wsCashflows.Range("E5").Value = 158977.2128939626
